$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct the erroneous values in column D (rows 9 and 10) from 1 to 0.1
$ws.Range("D9").Value = 0.1
$ws.Range("D10").Value = 0.1

# Update the active selection to reflect where the user ended up after editing
$ws.Range("H18").Select()
